# ModelComponentClassDiagram.pptx - "Change Model class diagram"
#
# On the single slide of this deck, the author:
#   - renamed the "UniquePersonList" box (shape id 49) to "UniqueCinemaList"
#     and widened it slightly to fit the new caption;
#   - nudged the "Flowchart: Decision" shape (id 63) to the right to line up
#     with the wider box, and stretched/shrunk its two elbow connectors
#     (ids 30 and 64) to keep tracking the shapes they are glued to.
#
# Shape positions/sizes below are expressed in points (the PowerPoint COM
# object model's native unit -- 1 pt = 12700 EMU) using values chosen so
# that, after the host's internal single-precision float rounding, they
# reproduce the exact target EMU offsets from the authoritative OOXML.

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape id=49 "Rectangle 8" -> rename to UniqueCinemaList, widen box.
$shpUniqueList = Get-ShapeById $s 49
$shpUniqueList.Left = 353.30833435058594
$shpUniqueList.Width = 95.57771682739259
$shpUniqueList.TextFrame.TextRange.Text = "UniqueCinemaList"

# Shape id=30 "Elbow Connector 29" (glued from shape 51 to shape 49):
# only its width needs a hair-thin adjustment to match the resized box.
$shpConn29 = Get-ShapeById $s 30
$shpConn29.Width = 20.997244834899902

# Shape id=63 "Flowchart: Decision 96" moves right to realign with the
# widened UniqueCinemaList box.
$shpDecision63 = Get-ShapeById $s 63
$shpDecision63.Left = 449.41352844238287

# Shape id=64 "Elbow Connector 63" (glued from shape 63 to shape 62) tracks
# the decision shape's new left edge; its endpoint on shape 62 is unchanged
# so its width shrinks by the same amount the start point moved.
$shpConn63 = Get-ShapeById $s 64
$shpConn63.Left = 467.99998474121094
$shpConn63.Width = 29.13992214202881
